$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A40").Value = "t1"
$ws.Range("A41").Value = "t2"
$ws.Range("A42").Value = "t3"

$i1 = $ws.Range("A40").Interior
$i1.TintAndShade = -0.499984740745262
Write-Host "after tint only:"

$i2 = $ws.Range("A41").Interior
$i2.ThemeColor = 2
Write-Host "after theme only:"
